$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2143.75
$ws.Range("I94").Value = 2143.75
$ws.Range("K94").Value = 2143.75
$ws.Range("M94").Value = -1692.75

$ws.Range("H98").Value = 8554.6
$ws.Range("J98").Value = 24832.666
$ws.Range("L98").Value = 24832.666
$ws.Range("N98").Value = -27828.666

$ws.Range("H113").Value = 3749.3333
$ws.Range("J113").Value = 4624.5
$ws.Range("L113").Value = 4624.5
$ws.Range("N113").Value = -11132.5

$ws.Range("H116").Value = 5499.6665
$ws.Range("I116").Value = 6499.5
$ws.Range("K116").Value = 6499.5
$ws.Range("M116").Value = -3057.5

$ws.Range("H122").Value = 8554.6
$ws.Range("J122").Value = 24832.666
$ws.Range("L122").Value = 74497.99800000001
$ws.Range("N122").Value = -79397.99800000001

$ws.Range("H137").Value = 38487.37
$ws.Range("I137").Value = 1575.8334
$ws.Range("J137").Value = 112310.445
$ws.Range("K137").Value = 4727.5002
$ws.Range("L137").Value = 336931.335
$ws.Range("M137").Value = -2177.5002
$ws.Range("N137").Value = -342031.335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 157.33333
$ws.Range("I5").Value = 133.92308
$ws.Range("K5").Value = 133.92308
$ws.Range("M5").Value = -21.92308

$ws.Range("H45").Value = 1004330.8
$ws.Range("I45").Value = 2004801.4
$ws.Range("J45").Value = 3860.2
$ws.Range("K45").Value = 2004801.4
$ws.Range("L45").Value = 3860.2
$ws.Range("M45").Value = -2004424.4
$ws.Range("N45").Value = -4614.2

$ws.Range("H61").Value = 1318.3
$ws.Range("I61").Value = 1131.4445
$ws.Range("K61").Value = 1131.4445
$ws.Range("M61").Value = -919.4445000000001

$ws.Range("H122").Value = 14705.5625
$ws.Range("I122").Value = 15486
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 46458
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -44008
$ws.Range("N122").Value = -13897

$ws.Range("H132").Value = 11318.857
$ws.Range("I132").Value = 12372.723
$ws.Range("K132").Value = 37118.169
$ws.Range("M132").Value = -34588.169

$ws.Range("H136").Value = 1318.3
$ws.Range("I136").Value = 1131.4445
$ws.Range("K136").Value = 3394.3335
$ws.Range("M136").Value = -844.3335000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 157.33333
$ws.Range("I4").Value = 133.92308
$ws.Range("K4").Value = 133.92308
$ws.Range("M4").Value = -18.92308

$ws.Range("H15").Value = 39499
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H29").Value = 2672
$ws.Range("I29").Value = 1508
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 1508
$ws.Range("L29").Value = 5000
$ws.Range("M29").Value = -1219
$ws.Range("N29").Value = -5578

$ws.Range("H138").Value = 55998.75
$ws.Range("J138").Value = 55998.75
$ws.Range("L138").Value = 55998.75
$ws.Range("N138").Value = -66278.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 313.73172
$ws.Range("I7").Value = 183.96
$ws.Range("K7").Value = 183.96
$ws.Range("M7").Value = -70.96000000000001

$ws.Range("H31").Value = 1258.4642
$ws.Range("I31").Value = 1258.4642
$ws.Range("K31").Value = 1258.4642
$ws.Range("M31").Value = -963.4641999999999

$ws.Range("H34").Value = 1258.4642
$ws.Range("I34").Value = 1258.4642
$ws.Range("K34").Value = 1258.4642
$ws.Range("M34").Value = -1056.4642

$ws.Range("H69").Value = 13194
$ws.Range("I69").Value = 8888
$ws.Range("K69").Value = 8888
$ws.Range("M69").Value = -8139

$ws.Range("H72").Value = 13194
$ws.Range("I72").Value = 8888
$ws.Range("K72").Value = 26664
$ws.Range("M72").Value = -22920

$ws.Range("H123").Value = 94650.42999999999
$ws.Range("J123").Value = 94650.42999999999
$ws.Range("L123").Value = 94650.42999999999
$ws.Range("N123").Value = -104450.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 250087.25
$ws.Range("I29").Value = 333416.34
$ws.Range("K29").Value = 1000249.02
$ws.Range("M29").Value = -999972.02

$ws.Range("H122").Value = 758.10254
$ws.Range("J122").Value = 775.96875
$ws.Range("L122").Value = 6983.71875
$ws.Range("N122").Value = -11883.71875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 16006160
$ws.Range("I3").Value = 5004900
$ws.Range("J3").Value = 23340332
$ws.Range("K3").Value = 5004900
$ws.Range("L3").Value = 23340332
$ws.Range("M3").Value = -5004784
$ws.Range("N3").Value = -23340564

$ws.Range("H10").Value = 6087
$ws.Range("I10").Value = 5375
$ws.Range("J10").Value = 7036.3335
$ws.Range("K10").Value = 5375
$ws.Range("L10").Value = 7036.3335
$ws.Range("M10").Value = -5206
$ws.Range("N10").Value = -7374.3335

$ws.Range("H11").Value = 50841200
$ws.Range("I11").Value = 24978500
$ws.Range("J11").Value = 80398570
$ws.Range("K11").Value = 24978500
$ws.Range("L11").Value = 80398570
$ws.Range("M11").Value = -24978361
$ws.Range("N11").Value = -80398848

$ws.Range("H21").Value = 27055.5
$ws.Range("I21").Value = 20500
$ws.Range("K21").Value = 20500
$ws.Range("M21").Value = -20327

$ws.Range("H24").Value = 33342800
$ws.Range("I24").Value = 100002500
$ws.Range("J24").Value = 12948.75
$ws.Range("K24").Value = 100002500
$ws.Range("L24").Value = 12948.75
$ws.Range("M24").Value = -100002327
$ws.Range("N24").Value = -13294.75

$ws.Range("H30").Value = 27055.5
$ws.Range("I30").Value = 20500
$ws.Range("K30").Value = 20500
$ws.Range("M30").Value = -20395

$ws.Range("H35").Value = 33044.2
$ws.Range("I35").Value = 31000
$ws.Range("J35").Value = 36110.5
$ws.Range("K35").Value = 31000
$ws.Range("L35").Value = 36110.5
$ws.Range("M35").Value = -30702
$ws.Range("N35").Value = -36706.5

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 3775
$ws.Range("I126").Value = 3400
$ws.Range("K126").Value = 10200
$ws.Range("M126").Value = -7730

$ws.Range("H132").Value = 3324
$ws.Range("I132").Value = 3399.389
$ws.Range("J132").Value = 3154.375
$ws.Range("K132").Value = 10198.167
$ws.Range("L132").Value = 9463.125
$ws.Range("M132").Value = -7668.167000000001
$ws.Range("N132").Value = -14523.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 48680.832
$ws.Range("J20").Value = 56737
$ws.Range("L20").Value = 56737
$ws.Range("N20").Value = -57189

$ws.Range("H42").Value = 28599.6
$ws.Range("I42").Value = 24500
$ws.Range("J42").Value = 31332.666
$ws.Range("K42").Value = 24500
$ws.Range("L42").Value = 31332.666
$ws.Range("M42").Value = -23937
$ws.Range("N42").Value = -32458.666

$ws.Range("H44").Value = 29976
$ws.Range("J44").Value = 29976
$ws.Range("L44").Value = 29976
$ws.Range("N44").Value = -30888

$ws.Range("H49").Value = 28599.6
$ws.Range("I49").Value = 24500
$ws.Range("J49").Value = 31332.666
$ws.Range("K49").Value = 24500
$ws.Range("L49").Value = 31332.666
$ws.Range("M49").Value = -24353
$ws.Range("N49").Value = -31626.666

$ws.Range("H56").Value = 19999.5
$ws.Range("J56").Value = 19999.5
$ws.Range("L56").Value = 19999.5
$ws.Range("N56").Value = -21381.5

$ws.Range("H136").Value = 5518.8335
$ws.Range("I136").Value = 4279.5
$ws.Range("K136").Value = 12838.5
$ws.Range("M136").Value = -10288.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 33493.332
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 33493.332
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 33493.332
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -33953.332

$ws.Range("H107").Value = 20833786
$ws.Range("J107").Value = 83333950
$ws.Range("L107").Value = 250001850
$ws.Range("N107").Value = -250005690

$ws.Range("H113").Value = 713.8182
$ws.Range("I113").Value = 385
$ws.Range("J113").Value = 810.5294
$ws.Range("K113").Value = 1155
$ws.Range("L113").Value = 2431.5882
$ws.Range("M113").Value = 1015
$ws.Range("N113").Value = -6771.5882

$ws.Range("H136").Value = 737.5161000000001
$ws.Range("I136").Value = 728.76666
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2186.29998
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = 363.7000200000002
$ws.Range("N136").Value = -8100
